# "Final amount received. Happy :)"
# The sheet tracked an expected PayPal payout in row 11 (Dollars/Rupees,
# formula-driven off row 10). Row 13 used to just hold a "Got in Bank"
# label in H13. The actual amount has now arrived, so:
#   - G13 gets the final received amount (in Rupees), formatted the same
#     way as the other Rupee cells in this block (copy format from G11).
#   - H13's label is replaced with a highlighted "<<< Got in Bank <<<"
#     marker (bold, green, centered) to call it out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Final amount received in the bank (Rupees), same number format/fill/
# border/font as the rest of the row-11 Rupee block.
$ws.Range("G13").Value = 229333.81
$ws.Range("G11").Copy()
$ws.Range("G13").PasteSpecial(-4122)  # xlPasteFormats

# Replace the old plain "Got in Bank" label with a highlighted marker.
$ws.Range("H13").Value = "<<< Got in Bank <<<"
$ws.Range("H13").Font.Bold = $true
$ws.Range("H13").Font.Color = 5287936   # RGB(0,176,80) -> OLE BGR 0x00B050
$ws.Range("H13").HorizontalAlignment = -4108  # xlCenter

# Move the active selection (as left by the editing session).
$null = $ws.Range("F22").Select()
